$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Final target values for the schedule grid (A1:F17).
$data = @{
    4  = @("8:40",  "-",         "João Rodrigues-Desenho Técnico", "-",         "Andre Lucca-Circuitos Elétricos", "-")
    6  = @("9:50",  "-",         "-",                               "-",         "-",                               "-")
    7  = @("10:40", "-",         "-",                               "-",         "João Rodrigues-Desenho Técnico",  "-")
    8  = @("11:30", "-",         "-",                               "-",         "-",                               "-")
    9  = @("12:20", "Almoço",    "Almoço",                          "Almoço",    "Almoço",                          "Almoço")
    10 = @("13:00", "-",         "-",                               "-",         "-",                               "-")
    11 = @("13:50", "-",         "-",                               "-",         "-",                               "-")
    12 = @("14:40", "-",         "-",                               "-",         "-",                               "-")
    13 = @("15:30", "Intervalo", "Intervalo",                       "Intervalo", "Intervalo",                       "Intervalo")
    14 = @("15:50", "-",         "-",                               "-",         "-",                               "-")
    15 = @("16:40", "-",         "-",                               "-",         "-",                               "-")
    16 = @("17:30", "-",         "-",                               "-",         "-",                               "-")
    17 = @("18:20", "",          "",                                "",          "",                                "")
}

foreach ($r in $data.Keys) {
    $rowVals = $data[$r]
    for ($c = 1; $c -le 6; $c++) {
        $ws.Cells.Item($r, $c).Value = $rowVals[$c - 1]
    }
}
